$d = $word.ActiveDocument

# --- Step 1: append the new sentence to the end of the
#     "It seems Mance relaxation requires ordered eigenvalues/eigenvectors.
#      This adds slightly to computational time. " paragraph.
#     We locate the insertion point via Find (rather than a cached Range),
#     collapse to its end, and InsertAfter there. The final sentence reads:
#     "...However, for Fortran this is only by 1ms for solid effect, and 9ms
#      for cross effect so overall negligible."
$fr = $d.Content
$found = $fr.Find.Execute("This adds slightly to computational time. ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$fr.Collapse(0)
$fr.InsertAfter("However, for Fortran this is only by 1ms for solid effect, and 9ms for cross effect so overall negligible.")

# --- Step 2: relocate the "_GoBack" bookmark (previously sitting alone in an
#     empty paragraph earlier in the document) to sit between "th" and "is is"
#     inside the sentence just inserted, i.e. splitting "this" -> "th" | "is".
#     Adding a bookmark whose name already exists elsewhere moves it (bookmark
#     names are unique per document), so the old empty-paragraph bookmark is
#     removed as a side effect. We re-resolve the split point fresh via Find
#     rather than reusing an old Range/position, since a Range's cached
#     Start/End can go stale across intervening mutations.
$fr2 = $d.Content
$found2 = $fr2.Find.Execute("However, for Fortran th", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$fr2.Collapse(0)
$d.Bookmarks.Add("_GoBack", $fr2)
